# Atualiza dados da planilha "vendas_atipicas":
#  - corrige alguns valores de linhas ja existentes
#  - insere uma nova linha (378216 / CABO HMASTON TIPO C IOS 30W) logo apos a
#    primeira linha (378212), deslocando as demais
#  - adiciona uma nova linha no final (386126 / CARREGADOR USB-C A GOLD 20W CA31-4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tabela final completa (linhas 2..10), na ordem em que devem aparecer.
$data = @(
    @("2025-07-07", 4, "BEMOL S/A",        "378212", 13546, "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON", -323, 1.1,  0.34),
    @("2025-07-07", 2, "BEMOL S/A",        "378216", 13538, "CABO HMASTON TIPO C IOS 30W",                                               -61, 1.1,  0.3),
    @("2025-07-09", 2, "BEMOL S/A",        "379513", 13000, "MOUSE PAD HARRY POTTER ESTAMPADO 26CMX21CM BLISTER C/1 UND LETRON",         -53, 1.03, 0.18),
    @("2025-07-09", 2, "MATHEUS SILVEIRA", "379106", 8915,  "CAPA IPHONE 11",                                                             28, 1.06, 0.24),
    @("2025-07-11", 2, "BEMOL S/A",        "380683", 14241, "MOCHILA PELUCIA STITCH",                                                    -27, 1.04, 0.19),
    @("2025-07-15", 2, "BEMOL S/A",        "383049", 12016, "PROJETOR ASTRONAUTA HMASTON",                                                -3, 1.04, 0.21),
    @("2025-07-16", 4, "BEMOL S/A",        "383665", 3984,  "BARALHO PLASTICO 1001 COPAG ESTOJO C/2 110 UNIDADES",                       -34, 1.22, 0.73),
    @("2025-07-17", 2, "BEMOL S/A",        "384275", 13185, "KIT LANCHE FUNDO DO MAR GARRAFA PLASTICA 500ML + MARMITA 700ML",            -34, 1.03, 0.17),
    @("2025-07-21", 2, "BEMOL S/A",        "386126", 10114, "CARREGADOR USB-C A GOLD 20W CA31-4",                                        -93, 1.05, 0.22)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
